# Update march 31, 2021 — append new Pakistan COVID-19 daily rows
# (report dates 2021-03-19 .. 2021-03-30, serial dates 44274..44285)
# to the bottom of the data table, replicating the row-389 formatting
# (date number format on B, colored "new cases" formulas on G/H/I).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each row: serial date, Confirmed, Deaths, Recovered, Active
$rows = @(
    @(44274, 623135, 13799, 579760, 29576),
    @(44275, 626802, 13843, 581852, 31107),
    @(44276, 630471, 13863, 583538, 33070),
    @(44277, 633741, 13935, 585271, 34535),
    @(44278, 637042, 13965, 586228, 36849),
    @(44279, 640988, 14028, 588975, 37985),
    @(44280, 645356, 14091, 591145, 40120),
    @(44281, 649824, 14158, 593282, 42384),
    @(44282, 654591, 14215, 595929, 44447),
    @(44283, 659116, 14256, 598197, 46663),
    @(44284, 663200, 14356, 600278, 48566),
    @(44285, 667957, 14434, 603126, 50397)
)

$startRow = 390
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $prev = $r - 1
    $data = $rows[$i]

    # Bring over number formats / fonts / fills from the row above so the
    # new row matches the existing "Country" table styling exactly.
    $ws.Range("A" + $prev + ":I" + $prev).Copy()
    $ws.Range("A" + $r + ":I" + $r).PasteSpecial(-4122)

    $ws.Range("A" + $r).Value = "Pakistan"
    $ws.Range("B" + $r).Value = $data[0]
    $ws.Range("C" + $r).Value = $data[1]
    $ws.Range("D" + $r).Value = $data[2]
    $ws.Range("E" + $r).Value = $data[3]
    $ws.Range("F" + $r).Value = $data[4]
    $ws.Range("G" + $r).Formula = "=C" + $r + "-C" + $prev
    $ws.Range("H" + $r).Formula = "=D" + $r + "-D" + $prev
    $ws.Range("I" + $r).Formula = "=E" + $r + "-E" + $prev
}

$excel.CutCopyMode = $false

# Scroll the view down to the new bottom rows and select the first empty
# row below the table, mirroring where Excel leaves the cursor after
# typing in the last data row.
$lastRow = $startRow + $rows.Count - 1
$excel.ActiveWindow.ScrollRow = 396
$ws.Range("A" + ($lastRow + 1)).Select()
